# Avances Controllers y DAOs.xlsx - progress update
# - Controllers sheet: scroll position moved (view-only, topLeftCell="A11")
# - Daos sheet: scroll position moved (topLeftCell="A83"), selection -> C92
# - Daos sheet: several "Completado" percentages bumped from 0 to 5%/10%
#   (new DAOs started but not finished yet)

$wb = $excel.ActiveWorkbook
$wsControllers = $wb.Worksheets.Item("Controllers")
$wsDaos = $wb.Worksheets.Item("Daos")

# --- Update progress values on the "Daos" sheet --------------------------
$wsDaos.Range("C2").Value = 0.05    # ActividadDAO
$wsDaos.Range("C48").Value = 0.05   # ObjetoDAO
$wsDaos.Range("C52").Value = 0.05   # PlanAdquisicionDAO
$wsDaos.Range("C53").Value = 0.05   # PlanAdquisicionPagoDAO
$wsDaos.Range("C59").Value = 0.05   # ProductoDAO
$wsDaos.Range("C91").Value = 0.1    # SubComponenteDAO
$wsDaos.Range("C96").Value = 0.05   # SubproductoDAO

# --- Window / selection state ---------------------------------------------
# Controllers sheet scrolled down, keeping its existing selection (C92)
$wsControllers.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1

# Daos sheet scrolled down and active cell moved to C92
$wsDaos.Activate()
$excel.ActiveWindow.ScrollRow = 83
$excel.ActiveWindow.ScrollColumn = 1
$wsDaos.Range("C92").Select()
